# Replace the floating "TODAY()" popup-style date with a fixed (blocking) date,
# per commit: "se reemplazo la ventana emergente con una ventana bloqueante info_procesos"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")
$ws.Activate()

# Inicio_del_proyecto (G3) used to be volatile TODAY(); pin it to a fixed date.
$ws.Range("G3").Formula = '=DATEVALUE("1/07/2022")'

# Task "1.1-Login y Tab Estadistica" progress updated to 50%.
$ws.Range("F9").Value = 0.5

# Restore the active cell/selection to where the user left off editing.
$ws.Range("W8").Select()
